$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.316.06"
$ws.Range("E2").Value = "  +0.25%  "

$ws.Range("D3").Value = "2.185.19"
$ws.Range("E3").Value = "  -1.31%  "

$ws.Range("E4").Value = "  +0.21%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "255.90"
$ws.Range("E5").Value = "  +5.48%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.627"
$ws.Range("E6").Value = "  +0.08%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "68.18"
$ws.Range("E7").Value = "  -2.47%  "

$ws.Range("E8").Value = "  +0.10%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.573"
$ws.Range("E9").Value = "  +3.67%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.36"
$ws.Range("E10").Value = "  -3.30%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "58.99"
$ws.Range("E11").Value = "  +2.00%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0936"
$ws.Range("E12").Value = "  -1.65%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.08"
$ws.Range("E13").Value = "  +5.36%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.105"
$ws.Range("E14").Value = "  +0.39%  "

$ws.Range("D15").Value = "2.516.33"
$ws.Range("E15").Value = "  -1.03%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.874"
$ws.Range("E16").Value = "  +4.16%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.44"
$ws.Range("E17").Value = "  -2.36%  "

$ws.Range("D18").Value = "2.166.71"
$ws.Range("E18").Value = "  -2.22%  "

$ws.Range("D19").Value = "41.276.91"
$ws.Range("E19").Value = "  +0.28%  "

$ws.Range("D20").Value = "0.0₃0958"
$ws.Range("E20").Value = "  +0.83%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.16"
$ws.Range("E21").Value = "  +1.07%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.95"
$ws.Range("E22").Value = "  -0.21%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "232.52"
$ws.Range("E23").Value = "  +0.37%  "

$ws.Range("E24").Value = "  +0.21%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.92"
$ws.Range("E25").Value = "  +8.99%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.77"
$ws.Range("E26").Value = "  +20.48%  "

$ws.Range("E27").Value = "  +0.01%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.51"

$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.17"
$ws.Range("E29").Value = "  -0.29%  "

$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "168.96"
$ws.Range("E30").Value = "  -2.05%  "

$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.62"
$ws.Range("E31").Value = "  +1.04%  "

$ws.Range("B32").Value = "Kaspa"
$ws.Range("C32").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.117"
$ws.Range("E32").Value = "  -1.73%  "

$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0756"
$ws.Range("E33").Value = "  +6.39%  "

$ws.Range("B34").Value = "Stellar"
$ws.Range("C34").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.123"
$ws.Range("E34").Value = "  -0.41%  "

$ws.Range("B35").Value = "InternetComputer(DFINITY)"
$ws.Range("C35").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.50"
$ws.Range("E35").Value = "  +5.31%  "

$ws.Range("B36").Value = "InjectiveProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "26.65"
$ws.Range("E36").Value = "  +12.01%  "

$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.18"
$ws.Range("E37").Value = "  +7.12%  "

$ws.Range("B38").Value = "Filecoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.61"
$ws.Range("E38").Value = "  +0.34%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0296"
$ws.Range("E39").Value = "  +6.14%  "

$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.20"
$ws.Range("E40").Value = "  -3.06%  "

$ws.Range("B41").Value = "Celestia"
$ws.Range("C41").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.37"
$ws.Range("E41").Value = "  +16.30%  "

$ws.Range("B42").Value = "THORChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.67"
$ws.Range("E42").Value = "  -2.75%  "

$ws.Range("B43").Value = "MultiversX"
$ws.Range("C43").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "63.66"
$ws.Range("E43").Value = "  -0.71%  "

$ws.Range("B44").Value = "FTXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.03"
$ws.Range("E44").Value = "  +0.83%  "

$ws.Range("B45").Value = "Algorand"
$ws.Range("C45").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.200"
$ws.Range("E45").Value = "  +1.28%  "

$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.64"
$ws.Range("E46").Value = "  -0.03%  "

$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.101"
$ws.Range("E47").Value = "  +1.56%  "

$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.19"
$ws.Range("E48").Value = "  +8.48%  "

$ws.Range("B49").Value = "BinanceUSD"
$ws.Range("C49").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.00"
$ws.Range("E49").Value = "  +0.20%  "

$ws.Range("B50").Value = "TrustWalletToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.18"
$ws.Range("E50").Value = "  +0.44%  "

$ws.Range("B51").Value = "SynthetixNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.25"
$ws.Range("E51").Value = "  -5.30%  "

